# closing dates were wrong for algo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F7"  = 544.6
    "G7"  = 561.65
    "H7"  = 542.25
    "I7"  = 556.25
    "J7"  = 548.2

    "G9"  = 547.9
    "H9"  = 541.65
    "I9"  = 544.2

    "G10" = 549.4
    "H10" = 543.65
    "I10" = 547.25

    "G11" = 550.8
    "H11" = 546.1
    "I11" = 548.8

    "G12" = 551.75
    "H12" = 546.3
    "I12" = 549.95

    "G13" = 550.2
    "H13" = 547
    "I13" = 547.15

    "G14" = 550.6
    "H14" = 547.15
    "I14" = 549.95

    "G15" = 552.6
    "H15" = 549.65
    "I15" = 552.1

    "G16" = 555
    "H16" = 551.6
    "I16" = 554.75

    "G17" = 556.9
    "H17" = 552.8
    "I17" = 554.5

    "G18" = 558.35
    "H18" = 553.55
    "I18" = 556.75

    "G19" = 559.65
    "H19" = 554
    "I19" = 559.3

    "G20" = 561.65
    "H20" = 557.25
    "I20" = 558.2

    "G21" = 559.45
    "H21" = 554.65
    "I21" = 557.2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
